$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Rayon_rotule (mm) : 8 -> 9.5
$ws.Range("B10").Value = 9.5

# Update Epaisseur_rotule (mm) : 6.75 -> 9
$ws.Range("B12").Value = 9

# Update active selection to B13 as left by the author when saving
$ws.Range("B13").Select()
